$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 170 (pushes existing rows 170..321 down to 172..323)
$ws.Rows.Item(170).Insert()
$ws.Rows.Item(171).Insert()

# --- New row 170 (Primera) ---
$ws.Cells.Item(170, 1).Value = 8
$ws.Cells.Item(170, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(170, 3).Value = "Coquimbo"
$ws.Cells.Item(170, 4).Value = 44789
$ws.Cells.Item(170, 5).Value = 4
$ws.Cells.Item(170, 6).Value = 100114014
$ws.Cells.Item(170, 7).Value = "Betarraga"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 2400
$ws.Cells.Item(170, 11).Value = 600
$ws.Cells.Item(170, 12).Value = 700
$ws.Cells.Item(170, 13).Value = 650
$ws.Cells.Item(170, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(170, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(170, 16).Value = 217
$ws.Cells.Item(170, 17).Value = 3
$ws.Cells.Item(170, 18).Value = "Hortaliza"

# --- New row 171 (Segunda) ---
$ws.Cells.Item(171, 1).Value = 8
$ws.Cells.Item(171, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(171, 3).Value = "Coquimbo"
$ws.Cells.Item(171, 4).Value = 44789
$ws.Cells.Item(171, 5).Value = 4
$ws.Cells.Item(171, 6).Value = 100114014
$ws.Cells.Item(171, 7).Value = "Betarraga"
$ws.Cells.Item(171, 8).Value = "Sin especificar"
$ws.Cells.Item(171, 9).Value = "Segunda"
$ws.Cells.Item(171, 10).Value = 1540
$ws.Cells.Item(171, 11).Value = 500
$ws.Cells.Item(171, 12).Value = 550
$ws.Cells.Item(171, 13).Value = 525
$ws.Cells.Item(171, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(171, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(171, 16).Value = 175
$ws.Cells.Item(171, 17).Value = 3
$ws.Cells.Item(171, 18).Value = "Hortaliza"
